# "Generate Report for Handoff" — update the localization-status report for
# the e2e\7bf67b17-951f-446d-bf93-370f1b21d917.md file: it has been
# re-handed-off, so its status moves from "Handed back: in sync with en-US"
# to "Ready for handoff", the handoff timestamps are refreshed, and an error
# detail noting the handback file is stale gets recorded per target language.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67dd622b3e240b35089460fe624bd64960f34477/e2e/7bf67b17-951f-446d-bf93-370f1b21d917.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5a9141b2df4552ec9a95d252391cc09c670445a/e2e/7bf67b17-951f-446d-bf93-370f1b21d917.md."

# --- "Overview" sheet: row 3 is the 7bf67b17-951f-446d-bf93-370f1b21d917.md file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus              # zh-cn column
$overview.Range("F3").Value = $newStatus              # de-de column
$overview.Range("G3").Value = "2016-10-18 03:23:47"   # Latest HO Xliff Generate Date

# --- "zh-cn" sheet: row 3 is the 7bf67b17-951f-446d-bf93-370f1b21d917.md file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus                  # Status
$zhcn.Range("H3").Value = "2016-10-18 03:23:25"       # Latest Handoff Datetime
$zhcn.Range("P3").Value = $errorDetail                # Error Detail
$zhcn.Columns.Item(16).ColumnWidth = 39.16666666666667  # widen Error Detail column to 40

# --- "de-de" sheet: row 3 is the 7bf67b17-951f-446d-bf93-370f1b21d917.md file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus                  # Status
$dede.Range("H3").Value = "2016-10-18 03:23:47"       # Latest Handoff Datetime
$dede.Range("P3").Value = $errorDetail                # Error Detail
$dede.Columns.Item(16).ColumnWidth = 39.16666666666667  # widen Error Detail column to 40
